$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '66.213.64'
$ws.Range("E2").Value = '  +3.29%  '
Set-TextValue "D3" '3.247.30'
$ws.Range("E3").Value = '  +7.28%  '
$ws.Range("E4").Value = '  +0.13%  '
Set-TextValue "D5" '582.21'
$ws.Range("E5").Value = '  +5.04%  '
Set-TextValue "D6" '152.48'
$ws.Range("E6").Value = '  +9.11%  '
$ws.Range("E7").Value = '  -0.22%  '
Set-TextValue "D8" '3.240.41'
$ws.Range("E8").Value = '  +7.41%  '
$ws.Range("E9").Value = '  +6.45%  '
Set-TextValue "D10" '7.09'
$ws.Range("E10").Value = '  +10.74%  '
$ws.Range("E11").Value = '  +7.39%  '
Set-TextValue "D12" '0.489'
$ws.Range("E12").Value = '  +6.18%  '
Set-TextValue "D13" '37.89'
$ws.Range("E13").Value = '  +4.86%  '
$ws.Range("E14").Value = '  +7.49%  '
Set-TextValue "D15" '3.770.81'
$ws.Range("E15").Value = '  +7.46%  '
Set-TextValue "D16" '66.279.26'
$ws.Range("E16").Value = '  +3.39%  '
Set-TextValue "D17" '552.48'
$ws.Range("E17").Value = '  +14.69%  '
Set-TextValue "D18" '3.251.53'
$ws.Range("E18").Value = '  +7.34%  '
$ws.Range("E19").Value = '  +3.03%  '
$ws.Range("E20").Value = '  +7.35%  '
Set-TextValue "D21" '14.56'
$ws.Range("E21").Value = '  +7.44%  '
Set-TextValue "D22" '0.746'
$ws.Range("E22").Value = '  +9.36%  '
Set-TextValue "D23" '7.89'
$ws.Range("E23").Value = '  +11.82%  '
Set-TextValue "D24" '13.52'
$ws.Range("E24").Value = '  +7.95%  '
Set-TextValue "D25" '81.48'
$ws.Range("E25").Value = '  +4.16%  '
$ws.Range("E26").Value = '  +0.12%  '
Set-TextValue "D27" '9.33'
$ws.Range("E27").Value = '  +19.82%  '
Set-TextValue "D28" '2.99'
$ws.Range("E28").Value = '  +10.25%  '
Set-TextValue "D29" '2.25'
$ws.Range("E29").Value = '  +7.48%  '
Set-TextValue "D30" '27.85'
$ws.Range("E30").Value = '  +8.25%  '
$ws.Range("E31").Value = '  +7.00%  '
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("E33").Value = '  +6.45%  '
Set-TextValue "D34" '565.85'
Set-TextValue "D35" '5.69'
$ws.Range("E35").Value = '  +5.37%  '
$ws.Range("E36").Value = '  +7.96%  '
Set-TextValue "D37" '55.24'
$ws.Range("E37").Value = '  +5.61%  '
$ws.Range("E38").Value = '  +13.57%  '
$ws.Range("E39").Value = '  +8.94%  '
$ws.Range("E40").Value = '  +7.82%  '
Set-TextValue "D41" '2.99'
$ws.Range("E41").Value = '  +10.72%  '
Set-TextValue "D42" '3.212.40'
$ws.Range("E42").Value = '  +11.74%  '
Set-TextValue "D43" '8.64'
$ws.Range("E43").Value = '  +4.64%  '
Set-TextValue "D44" '0.283'
$ws.Range("E44").Value = '  +17.05%  '
$ws.Range("E45").Value = '  +11.66%  '
Set-TextValue "D46" '26.57'
$ws.Range("E46").Value = '  +6.81%  '
$ws.Range("E47").Value = '  +0.07%  '
Set-TextValue "D48" '0.0₃0559'
$ws.Range("E48").Value = '  +5.91%  '
Set-TextValue "D49" '125.72'
$ws.Range("E49").Value = '  +4.91%  '
$ws.Range("E50").Value = '  +4.93%  '
$ws.Range("E51").Value = '  +9.85%  '
